$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string by writing the correct Arabic value for "Male" into B5.
# (Row 5 currently incorrectly shows the Arabic word for "Female"; code MLE means Male.)
$ws.Range("B5").Value = "الذكر"

# Autofit column B width (matches bestFit column width in target)
$ws.Columns("B").AutoFit() | Out-Null

# Set the active selection to D16 (as recorded in the workbook view)
$ws.Range("D16").Select() | Out-Null

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
